# added height in addSymbols
# Re-orders / updates the symbol rows (A2:F25) on Sheet1 to reflect the
# new "height" computation used when adding symbols.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1001, 18, 30, 75, 60, 72),
    @(501,  9,  52, 30, 75, 45),
    @(1201, 2,  10, 10, 10, 10),
    @(1203, 3,  15, 15, 15, 15),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6,  45, 30, 60, 45),
    @(1202, 2,  10, 10, 10, 10),
    @(401,  9,  48, 67, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(801,  3,  67, 65, 52, 45),
    @(101,  9,  30, 15, 60, 15),
    @(901,  16, 15, 45, 60, 60),
    @(601,  9,  60, 67, 60, 42),
    @(201,  9,  30, 15, 45, 30),
    @(1101, 0,  15, 30, 30, 0),
    @(2,    0,  2,  2,  2,  2),
    @(802,  0,  4,  5,  4,  0),
    @(1,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(3,    0,  3,  3,  3,  3),
    @(402,  0,  0,  4,  0,  0),
    @(602,  0,  0,  4,  0,  9),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
